$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Step 1: text/formatting edit, done as a precise surgical
# substitution on the WordprocessingML so the run split (and the new
# red-coloured "van ott internet" run) comes out exactly as in the
# target. We go through Content.WordOpenXML / InsertXML because this
# is a single run of text that needs to become four runs with
# different rPr, which Find/Replace can't express.
# ------------------------------------------------------------------
$xml = $d.Content.WordOpenXML

$oldPara = '<w:pPr><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">2 féle menüből is lehet választani, van ott internet, jól felszerelt, nagyon sok féle ízletes étel van, </w:t></w:r></w:p>'

$newPara = '<w:pPr><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">2 féle menüből is lehet választani, </w:t></w:r><w:r><w:rPr><w:color w:val="FF0000"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>van ott internet</w:t></w:r><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>, jól felszerelt, na</w:t></w:r><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>gyon sok féle ízletes étel van.</w:t></w:r></w:p>'

if ($xml.IndexOf($oldPara) -lt 0) {
    throw "menüből paragraph not found verbatim - aborting"
}
$xml = $xml.Replace($oldPara, $newPara)

# The _GoBack bookmark used to sit at the end of the "...Nagykanizsa."
# paragraph; it moves to the "menüből" paragraph below, so drop the old
# pair here (re-added further down via the Bookmarks API).
$oldTail = '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$newTail = '</w:p>'
if ($xml.IndexOf($oldTail) -lt 0) {
    throw "_GoBack bookmark tail not found verbatim - aborting"
}
$xml = $xml.Replace($oldTail, $newTail)

$d.Content.InsertXML($xml)

# ------------------------------------------------------------------
# Step 2: re-create the _GoBack bookmark around the edited paragraph,
# mirroring Word's own "remember the last edit spot" behaviour - start
# right at the top of the paragraph, end right after it.
# ------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*menüből*") {
        $d.Bookmarks.Add("_GoBack", $p.Range) | Out-Null
        break
    }
}
